$wb = $excel.ActiveWorkbook

# The table with the Input/Output "Source Name" headers lives on the
# "Events-Mulching" worksheet, inside the "annotationTable" ListObject.
$ws = $wb.Worksheets.Item("Events-Mulching")

# Rename the first column header: "Input [Source Name]" -> "Input [Sample Name]"
$ws.Range("A1").Value = "Input [Sample Name]"

# Rename the last column header: "Output [Source Name]" -> "Output [Sample Name]"
$ws.Range("AB1").Value = "Output [Sample Name]"
